$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.126.30"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "1.871.86"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "'0.5136"
$ws.Range("E7").Value = "  +1.95%  "

$ws.Range("E8").Value = "  +1.18%  "

$ws.Range("D9").Value = "'0.08367"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("D11").Value = "'41.67"
$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "'6.194"
$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.870.05"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "'7.292"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "'90.98"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").Value = "'0.06664"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "'6.033"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").Value = "28.165.20"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("D25").Value = "'2.249"
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").Value = "2.080.97"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").Value = "'2.474"
$ws.Range("E27").Value = "  -4.70%  "

$ws.Range("D28").Value = "'158.57"
$ws.Range("E28").Value = "  +1.39%  "

$ws.Range("D29").Value = "'20.56"
$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("D30").Value = "'124.86"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").Value = "'0.1060"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").Value = "'1.038"
$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("D33").Value = "'5.890"
$ws.Range("E33").Value = "  +4.39%  "

$ws.Range("D34").Value = "'3.603"
$ws.Range("E34").Value = "  -0.22%  "

$ws.Range("D35").Value = "'9.600"
$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("D36").Value = "'0.02438"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").Value = "'0.06537"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").Value = "'0.2186"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").Value = "'1.205"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("D40").Value = "'0.6496"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").Value = "'5.005"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("E43").Value = "  -0.47%  "

$ws.Range("D44").Value = "'0.6078"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("D45").Value = "'13.00"
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("D46").Value = "'3.679"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("E47").Value = "  -1.81%  "

$ws.Range("D48").Value = "'2.007"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").Value = "'1.216"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").Value = "'121.45"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "'0.06874"
$ws.Range("E51").Value = "  -0.59%  "
